$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.980.95"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.121.02"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.50"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.32"
$ws.Range("E6").Value = "  +1.37%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.115.10"
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("E10").Value = "  +8.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.78"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.464"
$ws.Range("E12").Value = "  -1.27%  "
$ws.Range("E13").Value = "  +2.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.26"
$ws.Range("E14").Value = "  +4.66%  "
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.633.49"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.825.41"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.15"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.123.23"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "464.62"
$ws.Range("E20").Value = "  +2.08%  "
$ws.Range("E21").Value = "  +1.43%  "
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.54"
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("E24").Value = "  -3.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.78"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.92"
$ws.Range("E27").Value = "  +7.18%  "
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("E33").Value = "  -3.25%  "
$ws.Range("E34").Value = "  +5.54%  "
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.39"
$ws.Range("E37").Value = "  +7.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.06"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.97"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "448.38"
$ws.Range("E40").Value = "  +4.35%  "
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.878.61"
$ws.Range("E43").Value = "  -3.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.276"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.78"
$ws.Range("E47").Value = "  +3.04%  "
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.53"
$ws.Range("E49").Value = "  -1.19%  "
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("E51").Value = "  -1.01%  "
